$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[3.6289936095376216, 9.264058714968186]"
$ws.Range("M2").Value = 0.00001031595251377304
$ws.Range("N2").Value = 0.00001031595251377304
$ws.Range("T2").Value = "[7.607071426249503, 11.13922609048369]"

# Row 3 updates
$ws.Range("L3").Value = "[4.01407389628571, 9.135283535936443]"
$ws.Range("M3").Value = 0.0000007749815289503914
$ws.Range("N3").Value = 0.000001549963057900783
$ws.Range("P3").Value = "[-3.383737432695313, -2.4780530640928498]"
$ws.Range("T3").Value = "[7.329565755787117, 10.401971636779647]"
$ws.Range("X3").Value = 9.777037037037209
$ws.Range("Y3").Value = 13.35037037037061
